# PlaneacionProyectoCotizadorAutos.xlsx - add US007 "Unit testing" user story
# + mark TA018(row18 of Iteracion 1)'s G-column progress as complete (G18: 0 -> 2)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Iteración 1" sheet: record progress on the existing TA018-ish task row
#    (G18 0 -> 2). The dependent SUM formulas on rows 17/5 recalc automatically.
# ---------------------------------------------------------------------------
$iter1 = $wb.Worksheets.Item("Iteración 1")
$iter1.Range("G18").Value = 2

# ---------------------------------------------------------------------------
# 2) Duplicate the US006 sheet as a starting point for the new US007 sheet,
#    keeping identical styling/merges/column widths.
# ---------------------------------------------------------------------------
$us006 = $wb.Worksheets.Item("US006")
$us006.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$us007 = $wb.Worksheets.Item($wb.Worksheets.Count)
$us007.Name = "US007"

# Header block values
$us007.Range("C3").Value = "US007"
$us007.Range("C4").Value = "Unit testing"
$us007.Range("C5").Value = "TA022"
$us007.Range("C6").Value = "Se desarrollará la tecnología necesaria para asegurar la integridad y correcta función del código desarrollado"

# Single task row (row 8) replaces the old 3-task rows copied from US006
$us007.Range("C8").Value = "TA022"
$us007.Range("D8").Value = "     Definición de tecnologia"
$us007.Range("E8").Value = "Se definirá la tecnologia a ser utilizada para este proposito"

# Copy the "last row" (thick bottom border) formatting from US001!C10:E10 onto
# US007 row 8, since it is now the single (and therefore last) task row.
$us006Fmt = $wb.Worksheets.Item("US001").Range("C10:E10")
$us006Fmt.Copy()
$us007.Range("C8:E8").PasteSpecial(-4122) | Out-Null

# Clear the two extra former task rows (9 & 10 from the US006 copy) content
$us007.Range("C9:E10").ClearContents()

# Row 11 (formerly the last-row styled task row) becomes a plain blank row:
# strip its thick-bottom / medium borders and number formatting down to the
# plain "blank row below table" look used elsewhere in the workbook.
$us007.Range("C11:E11").ClearContents()
$blankRows = $us007.Range("C9:E12")
$blankRows.Borders.LineStyle = -4142
$blankRows.VerticalAlignment = -4160
$us007.Range("C9:D11").NumberFormat = "@"
$us007.Range("E9:E11").WrapText = $true

$excel.Application.CutCopyMode = $false

$us007.Range("D12").Select()

Write-Host "Created US007 sheet"
